# Updates the crypto price/volume table (Price column D, Volume(1h) column E)
# with refreshed values from the latest scrape, matching the GitHub Actions
# commit "Updated cryptos list".
#
# Most Price-column values already contain two '.' separators (thousands +
# decimal), so Excel's type-inference leaves them as text automatically.
# A handful of Price values are plain decimals (e.g. "551.89") that Excel
# would otherwise auto-convert to a Number and round; for those we briefly
# capture/restore the cell's Style around a leading-apostrophe text literal
# so the value is forced back to text while leaving formatting untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.749.77"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "2.419.67"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.05%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.Value = "'551.89"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +0.01%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.Value = "'137.09"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("E9").Value = "  -0.08%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.Value = "'5.67"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -2.10%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.Value = "'0.353"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "2.850.47"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "59.717.01"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "2.414.41"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").Value = "  +2.38%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.Value = "'329.97"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -3.21%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.Value = "'1.00"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +0.07%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.Value = "'66.00"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +3.09%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +5.87%  "
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E29").Value = "  +0.13%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.Value = "'170.71"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -0.22%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.Value = "'18.66"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -2.16%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.Value = "'0.409"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -5.07%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.Value = "'312.56"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("E42").Value = "  +0.14%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.Value = "'138.61"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("E45").Value = "  +0.71%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.Value = "'19.50"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("E47").Value = "  +1.83%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.Value = "'0.402"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  +1.02%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.Value = "'17.56"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("E51").Value = "  -0.25%  "
